$d = $word.ActiveDocument

$pairs = @(
    @("282÷8=", "410÷2="),
    @("495÷9=", "168÷3="),
    @("567÷9=", "682÷6="),
    @("850÷6=", "956÷2="),
    @("995÷5=", "259÷3="),
    @("842÷9=", "414÷5="),
    @("370÷9=", "955÷5="),
    @("894÷8=", "333÷2="),
    @("692÷2=", "284÷9="),
    @("296÷7=", "585÷5="),
    @("253÷4=", "744÷3="),
    @("823÷2=", "510÷6="),
    @("852÷3=", "404÷4="),
    @("893÷8=", "555÷8="),
    @("264÷5=", "927÷7="),
    @("927÷8=", "893÷3="),
    @("562÷6=", "677÷4="),
    @("646÷7=", "207÷5="),
    @("142÷5=", "485÷6="),
    @("756÷9=", "340÷9="),
    @("716÷4=", "479÷3="),
    @("512÷9=", "932÷4="),
    @("228÷6=", "101÷9="),
    @("352÷9=", "377÷4="),
    @("269÷8=", "620÷3=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
